# Move VoSTR from fuels to io-model in acronym key variable list
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# --- capture the F6 formatting (cellXfs style 27) before it gets reassigned ---
# (F6 currently uses style index 27; the target state re-purposes that exact
# style slot for the relocated VoSTR row's F cell, while F6 itself moves to
# style index 6.)
$ws.Range("F6").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)  # xlPasteFormats -> stash style 27 formatting

# --- remove the VoSTR row from the "fuels" block (row 133) ---
$ws.Rows(133).Delete()

# --- insert a fresh row right before the "io-model" block's closing (thick-border) row ---
# After the delete above, the io-model block's final (thick-bottom-border) row is now row 182.
$ws.Rows(182).Insert()

# --- populate the newly inserted row 182 with the VoSTR content ---
$ws.Cells.Item(182, 1).Value = "io-model"
$ws.Cells.Item(182, 2).Value = "VoSTR"
$ws.Cells.Item(182, 3).Value = "VAT or Sales Tax Rate"
$ws.Cells.Item(182, 6).Value = "high"

# --- formatting: match the surrounding io-model rows (style index 2) for A:C ---
$ws.Range("A181").Copy()
$ws.Range("A182:C182").PasteSpecial(-4122)  # xlPasteFormats

# --- formatting: F182 reuses the stashed style-27 look (from the original F6) ---
$ws.Range("ZZ1").Copy()
$ws.Range("F182").PasteSpecial(-4122)  # xlPasteFormats

# --- clean up the scratch cell used to stash formatting ---
$ws.Range("ZZ1").ClearFormats()
$ws.Range("ZZ1").ClearContents()

# --- F6 itself switches off style 27 onto style 6's look (n/a -> still n/a, new fill) ---
$ws.Range("F40").Copy()
$ws.Range("F6").PasteSpecial(-4122)  # xlPasteFormats (F40 already uses style index 6)

$excel.CutCopyMode = 0

# --- view state: scroll back up to the top of the frozen pane, clear the old selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$ws.Range("A2").Select()
